$d = $word.ActiveDocument

# The template header used to show a "Date: {d[i].date}" line (right-
# justified) directly under the "CLIENT LOAN HISTORY REPORT" title.
# That whole paragraph (text + its paragraph mark) is removed, so the
# title is immediately followed by the blank right-justified paragraph
# that used to sit after the date line.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("{d[i].date}")) {
        $p.Range.Delete()
        break
    }
}
